{"js": "// Office.js (Word JavaScript API) script.\n// Applies the two edits described by the diff:\n//   1) In the paragraph \"Xi= -1,+1,-3,+3,-5,+5\" the three separate runs\n//      (\" -\", \"1,+\", \"1,-3,+3,-5,+5\") that were split around a couple of\n//      <w:proofErr/> grammar markers get merged back into one single run\n//      \" -1,+1,-3,+3,-5,+5\" (and the proofErr markers disappear).\n//   2) Near the end of the document, the first of a run of empty\n//      paragraphs (right after \"...almost no errors happen\") gets the\n//      text \"No value of S gives\", and the paragraph further down that\n//      holds \"14.16\" loses its <w:lastRenderedPageBreak/> marker.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst PKG_OPEN =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" ' +\n  'xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>';\nconst PKG_CLOSE =\n  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\n// Find the target paragraphs by their (unique) visible text rather than\n// hard-coded indexes, so the script is resilient to minor renumbering.\nlet mathPara = null;\nlet firstEmptyAfterHappen = null;\nlet para1416 = null;\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (t.indexOf(\"Xi=\") === 0 && t.indexOf(\"-1,+1,-3,+3,-5,+5\") !== -1) {\n    mathPara = items[i];\n  }\n  if (t.indexOf(\"almost no errors happen\") !== -1) {\n    // The next paragraph (first of the run of empty ones) is our target.\n    firstEmptyAfterHappen = items[i + 1];\n  }\n  if (t === \"14.16\") {\n    para1416 = items[i];\n  }\n}\n\n// --- Change 1: collapse the three math runs into a single run. ---\nif (mathPara) {\n  const ooxml = PKG_OPEN +\n    '<w:p w14:paraId=\"61864414\" w14:textId=\"3DA9D957\" w:rsidR=\"00B66724\" ' +\n    'w:rsidRDefault=\"00B66724\" w:rsidP=\"00920B42\">' +\n    '<w:pPr><w:ind w:firstLine=\"720\"/></w:pPr>' +\n    '<w:r><w:t>X</w:t></w:r>' +\n    '<w:r><w:rPr><w:vertAlign w:val=\"subscript\"/></w:rPr><w:t>i</w:t></w:r>' +\n    '<w:r><w:t>=</w:t></w:r>' +\n    '<w:r w:rsidR=\"00920B42\"><w:t xml:space=\"preserve\"> -1,+1,-3,+3,-5,+5</w:t></w:r>' +\n    '</w:p>' +\n    PKG_CLOSE;\n  mathPara.getRange(\"Whole\").insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// --- Change 2a: give the first empty paragraph the new sentence. ---\nif (firstEmptyAfterHappen) {\n  const ooxml = PKG_OPEN +\n    '<w:p w14:paraId=\"2B4EBBFE\" w14:textId=\"068ED645\" w:rsidR=\"00594BD5\" ' +\n    'w:rsidRDefault=\"00594BD5\" w:rsidP=\"00920B42\">' +\n    '<w:pPr><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr>' +\n    '<w:t>No value of S gives</w:t></w:r>' +\n    '</w:p>' +\n    PKG_CLOSE;\n  firstEmptyAfterHappen.getRange(\"Whole\").insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\n// --- Change 2b: drop the stale lastRenderedPageBreak marker. ---\nif (para1416) {\n  const ooxml = PKG_OPEN +\n    '<w:p w14:paraId=\"30BEC41F\" w14:textId=\"626D5B35\" w:rsidR=\"00594BD5\" ' +\n    'w:rsidRDefault=\"00594BD5\" w:rsidP=\"00920B42\">' +\n    '<w:pPr><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr></w:pPr>' +\n    '<w:r><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr>' +\n    '<w:t>14.16</w:t></w:r>' +\n    '</w:p>' +\n    PKG_CLOSE;\n  para1416.getRange(\"Whole\").insertOoxml(ooxml, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the two edits described by the diff:\n#   1) In the paragraph \"Xi= -1,+1,-3,+3,-5,+5\" the three separate runs\n#      (\" -\", \"1,+\", \"1,-3,+3,-5,+5\") that were split around a couple of\n#      <w:proofErr/> grammar markers get merged back into one single run\n#      \" -1,+1,-3,+3,-5,+5\" (and the proofErr markers disappear).\n#   2) Near the end of the document, the first of a run of empty\n#      paragraphs (right after \"...almost no errors happen\") gets the\n#      text \"No value of S gives\", and the paragraph further down that\n#      holds \"14.16\" loses its <w:lastRenderedPageBreak/> marker.\n\n$d = $word.ActiveDocument\n\n$mathParaIndex = $null\n$emptyAfterHappenIndex = $null\n$para1416Index = $null\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t.IndexOf(\"Xi=\") -eq 0 -and $t.IndexOf(\"-1,+1,-3,+3,-5,+5\") -ge 0) {\n        $mathParaIndex = $i\n    }\n    if ($t.IndexOf(\"almost no errors happen\") -ge 0) {\n        $emptyAfterHappenIndex = $i + 1\n    }\n    if ($t.Trim() -eq \"14.16\") {\n        $para1416Index = $i\n    }\n}\n\n# --- Change 1: collapse the three math runs into a single run. ---\nif ($mathParaIndex) {\n    $xml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body><w:p w14:paraId=\"61864414\" w14:textId=\"3DA9D957\" w:rsidR=\"00B66724\" w:rsidRDefault=\"00B66724\" w:rsidP=\"00920B42\"><w:pPr><w:ind w:firstLine=\"720\"/></w:pPr><w:r><w:t>X</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"subscript\"/></w:rPr><w:t>i</w:t></w:r><w:r><w:t>=</w:t></w:r><w:r w:rsidR=\"00920B42\"><w:t xml:space=\"preserve\"> -1,+1,-3,+3,-5,+5</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    [void]$d.Paragraphs.Item($mathParaIndex).Range.InsertXML($xml)\n}\n\n# --- Change 2a: give the first empty paragraph the new sentence. ---\nif ($emptyAfterHappenIndex) {\n    $xml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body><w:p w14:paraId=\"2B4EBBFE\" w14:textId=\"068ED645\" w:rsidR=\"00594BD5\" w:rsidRDefault=\"00594BD5\" w:rsidP=\"00920B42\"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr><w:t>No value of S gives</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    [void]$d.Paragraphs.Item($emptyAfterHappenIndex).Range.InsertXML($xml)\n}\n\n# --- Change 2b: drop the stale lastRenderedPageBreak marker. ---\nif ($para1416Index) {\n    $xml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\"><w:body><w:p w14:paraId=\"30BEC41F\" w14:textId=\"626D5B35\" w:rsidR=\"00594BD5\" w:rsidRDefault=\"00594BD5\" w:rsidP=\"00920B42\"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsiaTheme=\"minorEastAsia\"/></w:rPr><w:t>14.16</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    [void]$d.Paragraphs.Item($para1416Index).Range.InsertXML($xml)\n}\n"}
